$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 10067
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H54").Value = 25000
$ws.Range("J54").Value = 25000
$ws.Range("L54").Value = 25000
$ws.Range("N54").Value = -25972
$ws.Range("H87").Value = 127225.27
$ws.Range("J87").Value = 127225.27
$ws.Range("L87").Value = 127225.27
$ws.Range("N87").Value = -129721.27
$ws.Range("H90").Value = 127225.27
$ws.Range("J90").Value = 127225.27
$ws.Range("L90").Value = 381675.81
$ws.Range("N90").Value = -394155.81
$ws.Range("H98").Value = 383.41666
$ws.Range("I98").Value = 383.41666
$ws.Range("K98").Value = 383.41666
$ws.Range("M98").Value = 1114.58334
$ws.Range("H122").Value = 383.41666
$ws.Range("I122").Value = 383.41666
$ws.Range("K122").Value = 1150.24998
$ws.Range("M122").Value = 1299.75002
$ws.Range("H123").Value = 117362.75
$ws.Range("J123").Value = 117362.75
$ws.Range("L123").Value = 117362.75
$ws.Range("N123").Value = -127162.75
$ws.Range("H132").Value = 3573.484
$ws.Range("I132").Value = 3563.5
$ws.Range("K132").Value = 10690.5
$ws.Range("M132").Value = -8160.5
$ws.Range("H138").Value = 3952.7856
$ws.Range("J138").Value = 6554.2856
$ws.Range("L138").Value = 19662.8568
$ws.Range("N138").Value = -29942.8568
$ws.Range("H141").Value = 4615.5
$ws.Range("I141").Value = 4981.5
$ws.Range("J141").Value = 4249.5
$ws.Range("K141").Value = 14944.5
$ws.Range("L141").Value = 12748.5
$ws.Range("M141").Value = -9764.5
$ws.Range("N141").Value = -23108.5

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 44445
$ws.Range("J17").Value = 44445
$ws.Range("L17").Value = 44445
$ws.Range("N17").Value = -44791
$ws.Range("H45").Value = 1498.3611
$ws.Range("I45").Value = 1301.3636
$ws.Range("K45").Value = 1301.3636
$ws.Range("M45").Value = -924.3635999999999
$ws.Range("H46").Value = 18000
$ws.Range("J46").Value = 18000
$ws.Range("L46").Value = 18000
$ws.Range("N46").Value = -18638
$ws.Range("H122").Value = 3995.5
$ws.Range("I122").Value = 4012.7837
$ws.Range("J122").Value = 3782.3333
$ws.Range("K122").Value = 12038.3511
$ws.Range("L122").Value = 11346.9999
$ws.Range("M122").Value = -9588.3511
$ws.Range("N122").Value = -16246.9999
$ws.Range("H131").Value = 53995
$ws.Range("J131").Value = 53995
$ws.Range("L131").Value = 53995
$ws.Range("N131").Value = -64075
$ws.Range("H132").Value = 3466.0286
$ws.Range("I132").Value = 1733.0869
$ws.Range("K132").Value = 5199.2607
$ws.Range("M132").Value = -2669.2607

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 72500
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H94").Value = 760.1852
$ws.Range("I94").Value = 565.04
$ws.Range("K94").Value = 565.04
$ws.Range("M94").Value = -114.04
$ws.Range("H134").Value = 4485
$ws.Range("I134").Value = 1980.909
$ws.Range("K134").Value = 5942.727000000001
$ws.Range("M134").Value = -3407.727000000001

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 549.6053000000001
$ws.Range("I7").Value = 586.38464
$ws.Range("J7").Value = 469.91666
$ws.Range("K7").Value = 586.38464
$ws.Range("L7").Value = 469.91666
$ws.Range("M7").Value = -473.38464
$ws.Range("N7").Value = -695.91666
$ws.Range("H31").Value = 8554.772000000001
$ws.Range("I31").Value = 3759.5454
$ws.Range("K31").Value = 3759.5454
$ws.Range("M31").Value = -3464.5454
$ws.Range("H34").Value = 8554.772000000001
$ws.Range("I34").Value = 3759.5454
$ws.Range("K34").Value = 3759.5454
$ws.Range("M34").Value = -3557.5454
$ws.Range("H94").Value = 3736
$ws.Range("I94").Value = 3038.2
$ws.Range("J94").Value = 4123.6665
$ws.Range("K94").Value = 3038.2
$ws.Range("L94").Value = 4123.6665
$ws.Range("M94").Value = -2587.2
$ws.Range("N94").Value = -5025.6665
$ws.Range("H96").Value = 149562.38
$ws.Range("J96").Value = 149562.38
$ws.Range("L96").Value = 149562.38
$ws.Range("N96").Value = -155054.38
$ws.Range("H99").Value = 2310.1
$ws.Range("I99").Value = 1899.3334
$ws.Range("K99").Value = 1899.3334
$ws.Range("M99").Value = -401.3334
$ws.Range("H105").Value = 3060.75
$ws.Range("I105").Value = 3599.3333
$ws.Range("J105").Value = 1445
$ws.Range("K105").Value = 3599.3333
$ws.Range("L105").Value = 1445
$ws.Range("M105").Value = -1852.3333
$ws.Range("N105").Value = -4939
$ws.Range("H126").Value = 2310.1
$ws.Range("I126").Value = 1899.3334
$ws.Range("K126").Value = 5698.0002
$ws.Range("M126").Value = -3228.0002
$ws.Range("H134").Value = 5928.5
$ws.Range("I134").Value = 2970
$ws.Range("J134").Value = 12831.667
$ws.Range("K134").Value = 8910
$ws.Range("L134").Value = 38495.001
$ws.Range("M134").Value = -6375
$ws.Range("N134").Value = -43565.001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 428.15384
$ws.Range("I8").Value = 428.15384
$ws.Range("K8").Value = 1284.46152
$ws.Range("M8").Value = -1145.46152
$ws.Range("H113").Value = 1933
$ws.Range("J113").Value = 1933
$ws.Range("L113").Value = 5799
$ws.Range("N113").Value = -10139
$ws.Range("H129").Value = 15153096
$ws.Range("J129").Value = 41667956
$ws.Range("L129").Value = 125003868
$ws.Range("N129").Value = -125013868

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7898.9644
$ws.Range("I80").Value = 7557.2
$ws.Range("K80").Value = 7557.2
$ws.Range("M80").Value = -6559.2
$ws.Range("H83").Value = 7898.9644
$ws.Range("I83").Value = 7557.2
$ws.Range("K83").Value = 37786
$ws.Range("M83").Value = -32794
$ws.Range("H97").Value = 589.8823
$ws.Range("I97").Value = 682.5
$ws.Range("J97").Value = 457.57144
$ws.Range("K97").Value = 682.5
$ws.Range("L97").Value = 457.57144
$ws.Range("M97").Value = -186.5
$ws.Range("N97").Value = -1449.57144
$ws.Range("H122").Value = 4485.0835
$ws.Range("I122").Value = 3756.3215
$ws.Range("K122").Value = 11268.9645
$ws.Range("M122").Value = -8818.9645
$ws.Range("H126").Value = 4664.8335
$ws.Range("J126").Value = 10000
$ws.Range("L126").Value = 30000
$ws.Range("N126").Value = -34940
$ws.Range("H132").Value = 5627.6
$ws.Range("I132").Value = 5449.3213
$ws.Range("J132").Value = 6043.5835
$ws.Range("K132").Value = 16347.9639
$ws.Range("L132").Value = 18130.7505
$ws.Range("M132").Value = -13817.9639
$ws.Range("N132").Value = -23190.7505

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 51194.652
$ws.Range("I7").Value = 63097.11
$ws.Range("K7").Value = 63097.11
$ws.Range("M7").Value = -62985.11
$ws.Range("H14").Value = 6983.1665
$ws.Range("J14").Value = 7979.8
$ws.Range("L14").Value = 7979.8
$ws.Range("N14").Value = -8323.799999999999
$ws.Range("H16").Value = 1614.7142
$ws.Range("I16").Value = 1614.7142
$ws.Range("K16").Value = 1614.7142
$ws.Range("M16").Value = -1444.7142
$ws.Range("H40").Value = 3671.0557
$ws.Range("I40").Value = 2273
$ws.Range("K40").Value = 2273
$ws.Range("M40").Value = -2137
$ws.Range("H88").Value = 199500
$ws.Range("J88").Value = 199500
$ws.Range("L88").Value = 199500
$ws.Range("N88").Value = -200356
$ws.Range("H91").Value = 199500
$ws.Range("J91").Value = 199500
$ws.Range("L91").Value = 199500
$ws.Range("N91").Value = -202464
$ws.Range("H93").Value = 2050.0908
$ws.Range("I93").Value = 1867.2
$ws.Range("J93").Value = 2442
$ws.Range("K93").Value = 1867.2
$ws.Range("L93").Value = 2442
$ws.Range("M93").Value = -619.2
$ws.Range("N93").Value = -4938
$ws.Range("H122").Value = 5896.0625
$ws.Range("I122").Value = 6581.8335
$ws.Range("J122").Value = 3838.75
$ws.Range("K122").Value = 19745.5005
$ws.Range("L122").Value = 11516.25
$ws.Range("M122").Value = -17295.5005
$ws.Range("N122").Value = -16416.25
$ws.Range("H126").Value = 51194.652
$ws.Range("I126").Value = 63097.11
$ws.Range("K126").Value = 189291.33
$ws.Range("M126").Value = -186821.33
$ws.Range("H136").Value = 4017.6667
$ws.Range("I136").Value = 3669.6365
$ws.Range("J136").Value = 4974.75
$ws.Range("K136").Value = 11008.9095
$ws.Range("L136").Value = 14924.25
$ws.Range("M136").Value = -8458.9095
$ws.Range("N136").Value = -20024.25

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 352666.5
$ws.Range("I3").Value = 900000
$ws.Range("J3").Value = 243199.8
$ws.Range("K3").Value = 900000
$ws.Range("L3").Value = 243199.8
$ws.Range("M3").Value = -899886
$ws.Range("N3").Value = -243427.8
$ws.Range("H126").Value = 2554.111
$ws.Range("I126").Value = 2613.1428
$ws.Range("K126").Value = 7839.428400000001
$ws.Range("M126").Value = -5369.428400000001
$ws.Range("H132").Value = 5111.706
$ws.Range("I132").Value = 1850.7142
$ws.Range("K132").Value = 5552.142599999999
$ws.Range("M132").Value = -3022.142599999999
